# "D suite.xlsx" - update Results column values on the "Test Cases" sheet
# Runmode (C2:C5): "Y" -> "N"
# Results (D2:D5): "PASS" -> "SKIP"
# Selection moves from D2:D5 to C2:C5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

$ws.Range("C2:C5").Value = "N"
$ws.Range("D2:D5").Value = "SKIP"

$ws.Range("C2:C5").Select()
